# ex9.1.9 (Linear) Strong Stationary generator - alpha_zero
# "nuevos experimentos no convexos"
#
# The workbook stores every value (even the purely numeric-looking ones)
# as a shared string / text cell, not a native number. Plain
# `$range.Value = "0.93"` would let Excel's automatic type inference turn
# that into a real number cell, so numeric-looking replacements are
# written through a short dance that forces text storage and then
# restores the default ("Normal") style so no stray number-format/style
# is left behind on the cell.
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ---------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

# Row 2 (J_0_L0_v)
$ws.Range("A2").Value = "7.865 - x - 0.5y"
Set-TextValue $ws.Range("B2") "-5.865"
Set-TextValue $ws.Range("D2") "0.93"
Set-TextValue $ws.Range("E2") "4.3"
Set-TextValue $ws.Range("F2") "3.4000000000000004"

# Row 3 (J_0_L0_v)
$ws.Range("A3").Value = "-4.975 - 0.25x + y"
Set-TextValue $ws.Range("B3") "2.9749999999999996"
Set-TextValue $ws.Range("D3") "0.41"
Set-TextValue $ws.Range("E3") "4.699999999999999"
Set-TextValue $ws.Range("F3") "1.0"

# Row 4 (J_0_LP_v)
$ws.Range("A4").Value = "-7.865 + x + 0.5y"
Set-TextValue $ws.Range("B4") "-0.1349999999999998"
Set-TextValue $ws.Range("D4") "0.7"
Set-TextValue $ws.Range("E4") "8.5"
Set-TextValue $ws.Range("F4") "4.1"

# Row 5 (J_Ne_L0_v)
$ws.Range("A5").Value = "-11.559999999999999 + x - 2y"
Set-TextValue $ws.Range("B5") "-9.559999999999999"
Set-TextValue $ws.Range("D5") "0.36"
Set-TextValue $ws.Range("E5") "5.1"
Set-TextValue $ws.Range("F5") "2.5"

# Row 6 (J_Ne_L0_v)
$ws.Range("A6").Value = "-6.17 - y"
Set-TextValue $ws.Range("B6") "-6.17"
Set-TextValue $ws.Range("D6") "0.79"
Set-TextValue $ws.Range("E6") "9.5"
Set-TextValue $ws.Range("F6") "5.1"

# --- Punto_modificado ---------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws.Range("A2") "4.78"
Set-TextValue $ws.Range("B2") "6.17"

# --- Vector_bf -----------------------------------------------------------
# NOTE: Worksheets.Item(name) resolves case-insensitively, and "Vector_bf"
# / "Vector_BF" differ only by case, so both names would resolve to the
# same (first) sheet. Use the 1-based sheet index instead to target each
# one unambiguously (5 = Vector_bf, 6 = Vector_BF).
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "2.215"

# --- Vector_BF -----------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "-9.125"
Set-TextValue $ws.Range("A3") "11.9"

Write-Output "applied ex9.1.9 alpha_zero updates"
